$wb = $excel.ActiveWorkbook

# --- About sheet ---
$about = $wb.Worksheets.Item("About")

# Source year reference: 2017 -> 2015
$about.Range("B4").Value = 2015

# Title of the assessment report
$about.Range("B5").Value = "2015 Summer Reliability Assessment"

# Source URL
$about.Range("B6").Value = "http://www.nerc.com/pa/RAPA/ra/Reliability%20Assessments%20DL/2015_Summer_Reliability_Assessment.pdf"

# New citation detail row
$about.Range("B7").Value = "p.3, Table 1: Projected Demand, Resources, and Planning Reserve Margins, NERC Reference Margin Level (%)"

# Notes section: header stays "Notes"; body text is rewritten and now spans three rows
$about.Range("A9").Value = "Notes"
$about.Range("A10").Value = "The reserve margin (difference between the total generation available and the forecasted peak demand) in the U.S. "
$about.Range("A11").Value = "dataset doesn't vary by year, but the RM Reserve Margin variable is a time series to support countries that project "
$about.Range("A12").Value = "changes in future reserve margin by year."

# --- RM sheet ---
$rm = $wb.Worksheets.Item("RM")

# Units label, new cell
$rm.Range("A1").Value = "(dimensionless)"

# Updated reserve margin values across all years (B2:AK2)
$rm.Range("B2:AK2").Value = 0.1412
